$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 30000
$ws.Range("K82").Value = 90000
$ws.Range("M82").Value = -89594
$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 30000
$ws.Range("K85").Value = 90000
$ws.Range("M85").Value = -88596
$ws.Range("H125").Value = 26318894
$ws.Range("I125").Value = 45455910
$ws.Range("J125").Value = 5496.875
$ws.Range("K125").Value = 409103190
$ws.Range("L125").Value = 49471.875
$ws.Range("M125").Value = -409100730
$ws.Range("N125").Value = -54391.875
$ws.Range("H132").Value = 1343.1923
$ws.Range("I132").Value = 1362.42
$ws.Range("K132").Value = 4087.26
$ws.Range("M132").Value = -1557.26
$ws.Range("H138").Value = 3040272.2
$ws.Range("J138").Value = 3236149.5
$ws.Range("L138").Value = 9708448.5
$ws.Range("N138").Value = -9718728.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8995.925999999999
$ws.Range("I61").Value = 4438.353
$ws.Range("K61").Value = 4438.353
$ws.Range("M61").Value = -4226.353
$ws.Range("H64").Value = 42999.5
$ws.Range("J64").Value = 42999.5
$ws.Range("L64").Value = 42999.5
$ws.Range("N64").Value = -43495.5
$ws.Range("H67").Value = 42999.5
$ws.Range("J67").Value = 42999.5
$ws.Range("L67").Value = 42999.5
$ws.Range("N67").Value = -44715.5
$ws.Range("H74").Value = 37942.895
$ws.Range("J74").Value = 2776.8
$ws.Range("L74").Value = 2776.8
$ws.Range("N74").Value = -4524.8
$ws.Range("H77").Value = 37942.895
$ws.Range("J77").Value = 2776.8
$ws.Range("L77").Value = 13884
$ws.Range("N77").Value = -22620
$ws.Range("H122").Value = 2193.4524
$ws.Range("I122").Value = 2036.6923
$ws.Range("K122").Value = 6110.0769
$ws.Range("M122").Value = -3660.0769
$ws.Range("H132").Value = 2951564.2
$ws.Range("I132").Value = 5566291
$ws.Range("J132").Value = 9996.5
$ws.Range("K132").Value = 16698873
$ws.Range("L132").Value = 29989.5
$ws.Range("M132").Value = -16696343
$ws.Range("N132").Value = -35049.5
$ws.Range("H136").Value = 8995.925999999999
$ws.Range("I136").Value = 4438.353
$ws.Range("K136").Value = 13315.059
$ws.Range("M136").Value = -10765.059
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 50000476
$ws.Range("J80").Value = 503.875
$ws.Range("L80").Value = 503.875
$ws.Range("N80").Value = -2499.875
$ws.Range("H83").Value = 50000476
$ws.Range("J83").Value = 503.875
$ws.Range("L83").Value = 2519.375
$ws.Range("N83").Value = -12503.375
$ws.Range("H105").Value = 2887.2424
$ws.Range("I105").Value = 1842.76
$ws.Range("K105").Value = 1842.76
$ws.Range("M105").Value = -95.75999999999999
$ws.Range("H107").Value = 36291476
$ws.Range("I107").Value = 40179700
$ws.Range("K107").Value = 40179700
$ws.Range("M107").Value = -40177780
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -7830
$ws.Range("H134").Value = 5829.5884
$ws.Range("I134").Value = 2663.2
$ws.Range("K134").Value = 7989.599999999999
$ws.Range("M134").Value = -5454.599999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 6101.222
$ws.Range("I19").Value = 8659.166999999999
$ws.Range("K19").Value = 8659.166999999999
$ws.Range("M19").Value = -8489.166999999999
$ws.Range("H24").Value = 6101.222
$ws.Range("I24").Value = 8659.166999999999
$ws.Range("K24").Value = 8659.166999999999
$ws.Range("M24").Value = -8489.166999999999
$ws.Range("H31").Value = 4993.2144
$ws.Range("I31").Value = 1046.1666
$ws.Range("K31").Value = 1046.1666
$ws.Range("M31").Value = -751.1666
$ws.Range("H34").Value = 4993.2144
$ws.Range("I34").Value = 1046.1666
$ws.Range("K34").Value = 1046.1666
$ws.Range("M34").Value = -844.1666
$ws.Range("H58").Value = 6119.72
$ws.Range("I58").Value = 1392.2858
$ws.Range("J58").Value = 12136.454
$ws.Range("K58").Value = 1392.2858
$ws.Range("L58").Value = 12136.454
$ws.Range("M58").Value = -1189.2858
$ws.Range("N58").Value = -12542.454
$ws.Range("H69").Value = 25415.6
$ws.Range("I69").Value = 21359.666
$ws.Range("K69").Value = 21359.666
$ws.Range("M69").Value = -20610.666
$ws.Range("H72").Value = 25415.6
$ws.Range("I72").Value = 21359.666
$ws.Range("K72").Value = 64078.99800000001
$ws.Range("M72").Value = -60334.99800000001
$ws.Range("H93").Value = 13503.875
$ws.Range("I93").Value = 8764.666999999999
$ws.Range("J93").Value = 27721.5
$ws.Range("K93").Value = 8764.666999999999
$ws.Range("L93").Value = 27721.5
$ws.Range("M93").Value = -6892.666999999999
$ws.Range("N93").Value = -31465.5
$ws.Range("H107").Value = 1003.43243
$ws.Range("I107").Value = 309.95456
$ws.Range("J107").Value = 2020.5333
$ws.Range("K107").Value = 309.95456
$ws.Range("L107").Value = 2020.5333
$ws.Range("M107").Value = 1610.04544
$ws.Range("N107").Value = -5860.5333
$ws.Range("H134").Value = 6206.1113
$ws.Range("I134").Value = 1515.4615
$ws.Range("K134").Value = 4546.3845
$ws.Range("M134").Value = -2011.3845
$ws.Range("H136").Value = 6119.72
$ws.Range("I136").Value = 1392.2858
$ws.Range("J136").Value = 12136.454
$ws.Range("K136").Value = 4176.857400000001
$ws.Range("L136").Value = 36409.362
$ws.Range("M136").Value = -1626.857400000001
$ws.Range("N136").Value = -41509.362
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100653.35
$ws.Range("I2").Value = 154.64285
$ws.Range("J2").Value = 335150.34
$ws.Range("K2").Value = 927.8571000000001
$ws.Range("L2").Value = 2010902.04
$ws.Range("M2").Value = -814.8571000000001
$ws.Range("N2").Value = -2011128.04
$ws.Range("H3").Value = 1320
$ws.Range("I3").Value = 1320
$ws.Range("K3").Value = 3960
$ws.Range("M3").Value = -3848
$ws.Range("H24").Value = 3002.4167
$ws.Range("J24").Value = 2999.9092
$ws.Range("L24").Value = 8999.7276
$ws.Range("N24").Value = -9459.7276
$ws.Range("H42").Value = 13666.6
$ws.Range("J42").Value = 13666.6
$ws.Range("L42").Value = 40999.8
$ws.Range("N42").Value = -42067.8
$ws.Range("H107").Value = 1248.093
$ws.Range("J107").Value = 1713.826
$ws.Range("L107").Value = 5141.478
$ws.Range("N107").Value = -8981.477999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 47580000
$ws.Range("J11").Value = 62672500
$ws.Range("L11").Value = 62672500
$ws.Range("N11").Value = -62672778
$ws.Range("H57").Value = 79999.664
$ws.Range("J57").Value = 79999.664
$ws.Range("L57").Value = 79999.664
$ws.Range("N57").Value = -81639.664
$ws.Range("H126").Value = 5677.409
$ws.Range("I126").Value = 4146.769
$ws.Range("K126").Value = 12440.307
$ws.Range("M126").Value = -9970.307000000001
$ws.Range("H132").Value = 3398.7942
$ws.Range("I132").Value = 1859.6538
$ws.Range("J132").Value = 8401
$ws.Range("K132").Value = 5578.9614
$ws.Range("L132").Value = 25203
$ws.Range("M132").Value = -3048.9614
$ws.Range("N132").Value = -30263
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6909
$ws.Range("I7").Value = 5994.125
$ws.Range("K7").Value = 5994.125
$ws.Range("M7").Value = -5882.125
$ws.Range("H40").Value = 5301.933
$ws.Range("I40").Value = 4742.5
$ws.Range("K40").Value = 4742.5
$ws.Range("M40").Value = -4606.5
$ws.Range("H82").Value = 78965.92
$ws.Range("I82").Value = 100855.6
$ws.Range("K82").Value = 100855.6
$ws.Range("M82").Value = -100494.6
$ws.Range("H85").Value = 78965.92
$ws.Range("I85").Value = 100855.6
$ws.Range("K85").Value = 100855.6
$ws.Range("M85").Value = -99607.60000000001
$ws.Range("H126").Value = 6909
$ws.Range("I126").Value = 5994.125
$ws.Range("K126").Value = 17982.375
$ws.Range("M126").Value = -15512.375
$ws.Range("H132").Value = 5203.3335
$ws.Range("I132").Value = 3520.1282
$ws.Range("K132").Value = 10560.3846
$ws.Range("M132").Value = -8030.384600000001
$ws.Range("H136").Value = 9990.885
$ws.Range("I136").Value = 5633.3335
$ws.Range("J136").Value = 13725.929
$ws.Range("K136").Value = 16900.0005
$ws.Range("L136").Value = 41177.787
$ws.Range("M136").Value = -14350.0005
$ws.Range("N136").Value = -46277.787
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2966.5
$ws.Range("I96").Value = 2450
$ws.Range("J96").Value = 3999.5
$ws.Range("K96").Value = 2450
$ws.Range("L96").Value = 3999.5
$ws.Range("M96").Value = -1077
$ws.Range("N96").Value = -6745.5
$ws.Range("H132").Value = 14723.422
$ws.Range("I132").Value = 12798.223
$ws.Range("K132").Value = 38394.669
$ws.Range("M132").Value = -35864.669
